$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.561.43'
$ws.Range('E2').Value = '  +1.52%  '

$ws.Range('D3').Value = '1.908.51'
$ws.Range('E3').Value = '  +3.30%  '

$ws.Range('E4').Value = '  +0.44%  '

$ws.Range('D5').Value = "'247.10"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.56%  '

$ws.Range('D6').Value = "'0.632"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.62%  '

$ws.Range('E7').Value = '  +0.38%  '

$ws.Range('D8').Value = "'42.20"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.53%  '

$ws.Range('D9').Value = "'0.339"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.11%  '

$ws.Range('D10').Value = "'0.0706"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.75%  '

$ws.Range('D11').Value = "'0.0997"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('D12').Value = '2.182.83'
$ws.Range('E12').Value = '  +3.28%  '

$ws.Range('E13').Value = '  +9.22%  '

$ws.Range('D14').Value = '1.908.30'

$ws.Range('D15').Value = "'0.692"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.52%  '

$ws.Range('E16').Value = '  +3.62%  '

$ws.Range('D17').Value = '35.528.70'
$ws.Range('E17').Value = '  +1.49%  '

$ws.Range('D18').Value = "'72.16"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.05%  '

$ws.Range('D19').Value = '0.0₃0821'
$ws.Range('E19').Value = '  +3.54%  '

$ws.Range('D20').Value = "'243.80"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.24%  '

$ws.Range('D21').Value = "'12.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.84%  '

$ws.Range('D22').Value = "'4.88"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.74%  '

$ws.Range('E23').Value = '  +0.38%  '

$ws.Range('D24').Value = "'2.30"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.72%  '

$ws.Range('D25').Value = "'2.23"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +21.03%  '

$ws.Range('D26').Value = "'172.40"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.65%  '

$ws.Range('E27').Value = '  +8.95%  '

$ws.Range('D28').Value = "'18.01"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.39%  '

$ws.Range('D30').Value = "'0.971"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +27.69%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.0569"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.44%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'4.11"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.40%  '

$ws.Range('D33').Value = "'4.20"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.60%  '

$ws.Range('E34').Value = '  +0.39%  '

$ws.Range('E35').Value = '  +5.31%  '

$ws.Range('D36').Value = "'2.02"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.27%  '

$ws.Range('D37').Value = "'1.35"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.83%  '

$ws.Range('E38').Value = '  +3.29%  '

$ws.Range('D39').Value = "'0.0205"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.86%  '

$ws.Range('D40').Value = "'91.06"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = "'15.66"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.42%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.355.58'
$ws.Range('E42').Value = '  +0.66%  '

$ws.Range('D43').Value = "'49.73"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +43.67%  '

$ws.Range('D44').Value = "'0.0597"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.22%  '

$ws.Range('D45').Value = "'2.37"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.57%  '

$ws.Range('D46').Value = "'12.79"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.16%  '

$ws.Range('E47').Value = '  +0.86%  '

$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = "'2.77"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.44%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'6.68"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.96%  '

$ws.Range('D50').Value = '2.092.56'
$ws.Range('E50').Value = '  +3.41%  '

$ws.Range('D51').Value = "'0.0691"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.16%  '
